$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Total:" label to include the Paypal note
$ws.Range("A20").Value = "Total: (Includes +3% for Paypal)"

# Fill in the two newly-populated "Details" column values
$ws.Range("C3").Value = 3
$ws.Range("C7").Value = 1.5

# Update the grand-total formula to add the 3% Paypal fee
$ws.Range("F20").Formula = "=SUM(F2:F19)+3/100*(SUM(F2:F19))"

# Move the active selection to C6, matching the saved view state
$ws.Range("C6").Select()
